$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$casesQuery = "MATCH (c:case)`n    WHERE c.gender='MALE'`nOPTIONAL MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)`nOPTIONAL MATCH (f:file)-[*]->(c)`nRETURN DISTINCT`n    c.case_id AS ``Case ID``,`n     ct.clinical_trial_designation AS ``Trial Code``,`n     a.arm_id AS Arm,`n      a.arm_drug AS ``Arm Treatment``,`nc.disease AS Diagnosis,`n  c.gender AS Gender,`n    c.race AS Race,`n    c.ethnicity AS Ethnicity`n "
$statQuery = "MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)`n    WHERE c.gender = `"MALE`"`nOPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)`nRETURN `n     COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,`n     COUNT(DISTINCT c.case_id) AS Cases,`n      COUNT(DISTINCT f) AS Files"
$filesQuery = "`nMATCH (c:case)`n WHERE c.gender='MALE'`n   MATCH (f:file)`n      OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)`n       MATCH (f)-[*]->(c)`n      OPTIONAL MATCH (f)-->(parent)`n      WITH`n        f, parent, c, a, ct, `n        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,`n        toInteger(floor(log(f.file_size)/log(1024))) as i,`n        2 as precision`nWITH`n        f, parent, c, a, ct,`n        f.file_size /(1024^i) AS value, `n        10^precision AS factor,`n        units[i] as unit`nWITH    `n        f, parent, c, a, ct, unit,`n        round(factor * value)/factor AS size`n      RETURN DISTINCT `n       f.file_name AS ``File Name``,`n       head(labels(parent)) as Association,`n       f.file_description AS Description,`n       f.file_format AS ``File Format``,`n         CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,`n       ct.clinical_trial_designation AS ``Trial Code``,`n       a.arm_id AS Arm,`n       c.case_id AS ``Case ID```n        "

$ws.Range("B2").Value = $casesQuery
$ws.Range("C2").Value = $statQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("B3").Value = $filesQuery

$ws.Rows.Item(2).RowHeight = 210
$ws.Rows.Item(3).RowHeight = 409.5

$ws.Columns.Item(1).ColumnWidth = 8.85546875
$ws.Range("B1:C1").ColumnWidth = 75.85546875
$ws.Columns.Item(4).ColumnWidth = 70.28515625
$ws.Columns.Item(5).ColumnWidth = 28.5703125

$ws.Range("D3").Select() | Out-Null

Write-Host "done"
